# Add an EXTRACT column to the find_and_replace sheet and make that sheet
# the active tab, with E2 selected as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("find_and_replace")

# Add the new header value, matching the text formatting of the other
# header cells in row 1.
$ws.Range("D1").Value = "EXTRACT"
$ws.Range("D1").NumberFormat = $ws.Range("C1").NumberFormat

# Activate this sheet, so it becomes the selected/active tab, and select E2.
$ws.Activate()
$ws.Range("E2").Select()
